# Refined metadata to be additional tab
#
# This workbook was re-pulled from the PanelApp API (hence the refreshed
# "time_taken" / query timestamps on the `data` sheet), and a new
# `metadata` sheet describing the panel/query itself was appended after it.

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# ---------------------------------------------------------------------
# 1) Refresh the per-row query timestamps on the `data` sheet (column F)
# ---------------------------------------------------------------------
$newTimes = @(
    "2021-10-05 14:22:19.551011",
    "2021-10-05 14:22:19.551018",
    "2021-10-05 14:22:19.551020",
    "2021-10-05 14:22:19.551022",
    "2021-10-05 14:22:19.551025",
    "2021-10-05 14:22:19.551027",
    "2021-10-05 14:22:19.551029",
    "2021-10-05 14:22:19.551031",
    "2021-10-05 14:22:19.551033",
    "2021-10-05 14:22:19.551035",
    "2021-10-05 14:22:19.551037",
    "2021-10-05 14:22:19.551039",
    "2021-10-05 14:22:19.551041",
    "2021-10-05 14:22:19.551043",
    "2021-10-05 14:22:19.551045",
    "2021-10-05 14:22:19.551047",
    "2021-10-05 14:22:19.551049",
    "2021-10-05 14:22:19.551051",
    "2021-10-05 14:22:19.551053",
    "2021-10-05 14:22:19.551055",
    "2021-10-05 14:22:19.551057",
    "2021-10-05 14:22:19.551060",
    "2021-10-05 14:22:19.551061",
    "2021-10-05 14:22:19.551063",
    "2021-10-05 14:22:19.551066",
    "2021-10-05 14:22:19.551068",
    "2021-10-05 14:22:19.551070",
    "2021-10-05 14:22:19.551072",
    "2021-10-05 14:22:19.551074",
    "2021-10-05 14:22:19.551076",
    "2021-10-05 14:22:19.551078",
    "2021-10-05 14:22:19.551081",
    "2021-10-05 14:22:19.551084",
    "2021-10-05 14:22:19.551086",
    "2021-10-05 14:22:19.551088",
    "2021-10-05 14:22:19.551090",
    "2021-10-05 14:22:19.551092",
    "2021-10-05 14:22:19.551094",
    "2021-10-05 14:22:19.551096",
    "2021-10-05 14:22:19.551098",
    "2021-10-05 14:22:19.551100",
    "2021-10-05 14:22:19.551102",
    "2021-10-05 14:22:19.551105",
    "2021-10-05 14:22:19.551107",
    "2021-10-05 14:22:19.551109",
    "2021-10-05 14:22:19.551111",
    "2021-10-05 14:22:19.551113",
    "2021-10-05 14:22:19.551115",
    "2021-10-05 14:22:19.551117",
    "2021-10-05 14:22:19.551119",
    "2021-10-05 14:22:19.551121",
    "2021-10-05 14:22:19.551123",
    "2021-10-05 14:22:19.551126",
    "2021-10-05 14:22:19.551128",
    "2021-10-05 14:22:19.551130",
    "2021-10-05 14:22:19.551132"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $newTimes[$i]
}

# ---------------------------------------------------------------------
# 2) Add the new `metadata` sheet right after `data`
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Add($null, $dataSheet)
$meta.Name = "metadata"

$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Primary lymphoedema"
$meta.Range("C2").Value = 65

# data_version "2.18" must stay text, not become the number 2.18
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "2.18"

$meta.Range("E2").Value = "2021-09-08T13:47:13.257858Z"
$meta.Range("F2").Value = "2021-10-05 14:22:19.548613"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/65/?format=json"

# Header row (B1:G1) + the index cell (A2) use the same bold/centered/
# bordered style as the `data` sheet's header row / index column.
$headerRange = $meta.Range("B1:G1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$indexCell = $meta.Range("A2")
$indexCell.Font.Bold = $true
$indexCell.HorizontalAlignment = -4108
$indexCell.VerticalAlignment = -4160
$indexCell.Borders.LineStyle = 1

# Put the selection back on the `data` sheet / A1, matching the original.
$dataSheet.Activate()
$dataSheet.Range("A1").Select() | Out-Null
